$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
  @(4, "CC.EST"),
  @(5, "CC.NO.SRC"),
  @(6, "CC.PER.RNK"),
  @(7, "CC.PER.RNK.LOWER"),
  @(8, "CC.PER.RNK.UPPER"),
  @(9, "CC.STD.ERR"),
  @(10, "DT.ODA.ALLD.CD"),
  @(11, "DT.ODA.ALLD.KD"),
  @(12, "DT.ODA.OATL.CD"),
  @(13, "DT.ODA.OATL.KD"),
  @(14, "DT.ODA.ODAT.CD"),
  @(15, "DT.ODA.ODAT.GI.ZS"),
  @(16, "DT.ODA.ODAT.GN.ZS"),
  @(17, "DT.ODA.ODAT.KD"),
  @(18, "DT.ODA.ODAT.MP.ZS"),
  @(19, "DT.ODA.ODAT.PC.ZS"),
  @(20, "DT.ODA.ODAT.XP.ZS"),
  @(21, "GE.EST"),
  @(22, "GE.NO.SRC"),
  @(23, "GE.PER.RNK"),
  @(24, "GE.PER.RNK.LOWER"),
  @(25, "GE.PER.RNK.UPPER"),
  @(26, "GE.STD.ERR"),
  @(27, "PV.EST"),
  @(28, "PV.NO.SRC"),
  @(29, "PV.PER.RNK"),
  @(30, "PV.PER.RNK.LOWER"),
  @(31, "PV.PER.RNK.UPPER"),
  @(32, "PV.STD.ERR"),
  @(33, "RQ.EST"),
  @(34, "RQ.NO.SRC"),
  @(35, "RQ.PER.RNK"),
  @(36, "RQ.PER.RNK.LOWER"),
  @(37, "RQ.PER.RNK.UPPER"),
  @(38, "RQ.STD.ERR"),
  @(39, "RL.ES")
)

foreach ($row in $aValues) {
  $r = $row[0]
  $a = $row[1]
  $ws.Cells.Item($r, 1).Value = $a
}

$bValues = @(
  @(10, "Net official development assistance and official aid received (current US`$)"),
  @(11, "Net official development assistance and official aid received (constant 2021 US`$)"),
  @(12, "Net official aid received (current US`$)"),
  @(13, "Net official aid received (constant 2021 US`$)"),
  @(14, "Net official development assistance received (current US`$)"),
  @(15, "Net ODA received (% of gross capital formation)"),
  @(16, "Net ODA received (% of GNI)"),
  @(17, "Net official development assistance received (constant 2021 US`$)"),
  @(18, "Net ODA received (% of imports of goods, services and primary income)"),
  @(19, "Net ODA received per capita (current US`$)"),
  @(20, "Net ODA received (% of central government expense)")
)

foreach ($row in $bValues) {
  $r = $row[0]
  $b = $row[1]
  $ws.Cells.Item($r, 2).Value = $b
}

# Column widths (target bestFit widths from the authored workbook; engine quantizes
# ColumnWidth to 1/6-character steps, so we pick the input that lands closest to the
# authored width after that internal rounding)
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws.Columns.Item(2).ColumnWidth = 64.66666666666667
$ws.Columns.Item(3).ColumnWidth = 11.0

# B21 ends up with its own (empty) cell style in the authored file - reproduce
# that by touching the font (forces a dedicated style record) without changing
# its visible appearance.
$ws.Range("B21").Font.ThemeColor = 1

# Selection change
$ws.Range("B21").Select()

# Page setup orientation
$ws.PageSetup.Orientation = 1
